# Auto-generated script: applies scheduled-runner market-price refresh
# to the Cactuar_Profits-style sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Every changed cell is a cached/static value (no formulas in this workbook),
# so each update is a direct Range.Value assignment.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19: H19,J19,L19,N19
$ws.Range("H19").Value = 1802.75
$ws.Range("J19").Value = 999
$ws.Range("L19").Value = 999
$ws.Range("N19").Value = -1349
# Row 64: H64,I64,K64,M64
$ws.Range("H64").Value = 14432101
$ws.Range("I64").Value = 4419812
$ws.Range("K64").Value = 4419812
$ws.Range("M64").Value = -4419564
# Row 67: H67,I67,K67,M67
$ws.Range("H67").Value = 14432101
$ws.Range("I67").Value = 4419812
$ws.Range("K67").Value = 4419812
$ws.Range("M67").Value = -4418954
# Row 74: H74,I74,K74,M74
$ws.Range("H74").Value = 5799.3335
$ws.Range("I74").Value = 4449.5
$ws.Range("K74").Value = 4449.5
$ws.Range("M74").Value = -3513.5
# Row 77: H77,I77,K77,M77
$ws.Range("H77").Value = 5799.3335
$ws.Range("I77").Value = 4449.5
$ws.Range("K77").Value = 22247.5
$ws.Range("M77").Value = -17567.5
# Row 107: H107,I107,K107,M107
$ws.Range("H107").Value = 541.2308
$ws.Range("I107").Value = 501.45456
$ws.Range("K107").Value = 501.45456
$ws.Range("M107").Value = 1418.54544
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 9078.835999999999
$ws.Range("I132").Value = 2505.7097
$ws.Range("K132").Value = 7517.1291
$ws.Range("M132").Value = -4987.1291
# Row 137: H137,I137,K137,M137
$ws.Range("H137").Value = 15157485
$ws.Range("I137").Value = 1499.5
$ws.Range("K137").Value = 4498.5
$ws.Range("M137").Value = -1948.5
# Row 138: H138,I138,K138,M138
$ws.Range("H138").Value = 4499.6064
$ws.Range("I138").Value = 912.05884
$ws.Range("K138").Value = 2736.17652
$ws.Range("M138").Value = 2403.82348

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2: H2,I2,K2,M2
$ws.Range("H2").Value = 2494882.2
$ws.Range("I2").Value = 4364133
$ws.Range("K2").Value = 4364133
$ws.Range("M2").Value = -4364020
# Row 32: H32,I32,K32,M32
$ws.Range("H32").Value = 5204.4526
$ws.Range("I32").Value = 2652.5
$ws.Range("K32").Value = 2652.5
$ws.Range("M32").Value = -2365.5
# Row 45: H45,I45,K45,M45
$ws.Range("H45").Value = 3158.842
$ws.Range("I45").Value = 2901.5334
$ws.Range("K45").Value = 2901.5334
$ws.Range("M45").Value = -2524.5334
# Row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 39752
$ws.Range("I61").Value = 51503
$ws.Range("J61").Value = 4499
$ws.Range("K61").Value = 51503
$ws.Range("L61").Value = 4499
$ws.Range("M61").Value = -51291
$ws.Range("N61").Value = -4923
# Row 74: H74,I74,K74,M74
$ws.Range("H74").Value = 19233054
$ws.Range("I74").Value = 35715844
$ws.Range("K74").Value = 35715844
$ws.Range("M74").Value = -35714970
# Row 77: H77,I77,K77,M77
$ws.Range("H77").Value = 19233054
$ws.Range("I77").Value = 35715844
$ws.Range("K77").Value = 178579220
$ws.Range("M77").Value = -178574852
# Row 116: H116,I116,K116,M116
$ws.Range("H116").Value = 2494882.2
$ws.Range("I116").Value = 4364133
$ws.Range("K116").Value = 4364133
$ws.Range("M116").Value = -4361839
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 20989.572
$ws.Range("I132").Value = 26171.088
$ws.Range("K132").Value = 78513.264
$ws.Range("M132").Value = -75983.264
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 39752
$ws.Range("I136").Value = 51503
$ws.Range("J136").Value = 4499
$ws.Range("K136").Value = 154509
$ws.Range("L136").Value = 13497
$ws.Range("M136").Value = -151959
$ws.Range("N136").Value = -18597

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3: H3,I3,K3,M3
$ws.Range("H3").Value = 2494882.2
$ws.Range("I3").Value = 4364133
$ws.Range("K3").Value = 4364133
$ws.Range("M3").Value = -4364019
# Row 20: H20,I20,J20,K20,L20,M20,N20
$ws.Range("H20").Value = 2707.4333
$ws.Range("I20").Value = 2254.261
$ws.Range("J20").Value = 4196.4287
$ws.Range("K20").Value = 2254.261
$ws.Range("L20").Value = 4196.4287
$ws.Range("M20").Value = -2007.261
$ws.Range("N20").Value = -4690.4287
# Row 132: H132,J132,L132,N132
$ws.Range("H132").Value = 41617.65
$ws.Range("J132").Value = 41617.65
$ws.Range("L132").Value = 41617.65
$ws.Range("N132").Value = -51737.65
# Row 134: H134,I134,K134,M134
$ws.Range("H134").Value = 4687.2856
$ws.Range("I134").Value = 1716.625
$ws.Range("K134").Value = 5149.875
$ws.Range("M134").Value = -2614.875
# Row 135: H135,J135,L135,N135
$ws.Range("H135").Value = 42751.035
$ws.Range("J135").Value = 42751.035
$ws.Range("L135").Value = 42751.035
$ws.Range("N135").Value = -52891.035

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 5: H5,I5,J5,K5,L5,M5,N5
$ws.Range("H5").Value = 1415.4
$ws.Range("I5").Value = 1856.3334
$ws.Range("J5").Value = 754
$ws.Range("K5").Value = 1856.3334
$ws.Range("L5").Value = 754
$ws.Range("M5").Value = -1744.3334
$ws.Range("N5").Value = -978
# Row 16: H16,I16,J16,K16,L16,M16,N16
$ws.Range("H16").Value = 1750.25
$ws.Range("I16").Value = 1786
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1786
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1499
$ws.Range("N16").Value = -2074
# Row 31: H31,I31,J31,K31,L31,M31,N31
$ws.Range("H31").Value = 4743
$ws.Range("I31").Value = 1692.619
$ws.Range("J31").Value = 9013.532999999999
$ws.Range("K31").Value = 1692.619
$ws.Range("L31").Value = 9013.532999999999
$ws.Range("M31").Value = -1397.619
$ws.Range("N31").Value = -9603.532999999999
# Row 34: H34,I34,J34,K34,L34,M34,N34
$ws.Range("H34").Value = 4743
$ws.Range("I34").Value = 1692.619
$ws.Range("J34").Value = 9013.532999999999
$ws.Range("K34").Value = 1692.619
$ws.Range("L34").Value = 9013.532999999999
$ws.Range("M34").Value = -1490.619
$ws.Range("N34").Value = -9417.532999999999
# Row 52: H52,J52,L52,N52
$ws.Range("H52").Value = 53390
$ws.Range("J52").Value = 52780
$ws.Range("L52").Value = 52780
$ws.Range("N52").Value = -53368
# Row 58: H58,I58,J58,K58,L58,M58,N58
$ws.Range("H58").Value = 773330.3
$ws.Range("I58").Value = 913072.2
$ws.Range("J58").Value = 4750
$ws.Range("K58").Value = 913072.2
$ws.Range("L58").Value = 4750
$ws.Range("M58").Value = -912869.2
$ws.Range("N58").Value = -5156
# Row 62: H62,I62,J62,K62,L62,M62,N62
$ws.Range("H62").Value = 40764.375
$ws.Range("I62").Value = 2004.5
$ws.Range("J62").Value = 53684.332
$ws.Range("K62").Value = 2004.5
$ws.Range("L62").Value = 53684.332
$ws.Range("M62").Value = -1380.5
$ws.Range("N62").Value = -54932.332
# Row 65: H65,I65,J65,K65,L65,M65,N65
$ws.Range("H65").Value = 40764.375
$ws.Range("I65").Value = 2004.5
$ws.Range("J65").Value = 53684.332
$ws.Range("K65").Value = 10022.5
$ws.Range("L65").Value = 268421.66
$ws.Range("M65").Value = -6902.5
$ws.Range("N65").Value = -274661.66
# Row 113: H113,I113,J113,K113,L113,M113,N113
$ws.Range("H113").Value = 1750.25
$ws.Range("I113").Value = 1786
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1786
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 384
$ws.Range("N113").Value = -5840
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 18538620
$ws.Range("I132").Value = 20848822
$ws.Range("K132").Value = 62546466
$ws.Range("M132").Value = -62543936
# Row 134: H134,J134,L134,N134
$ws.Range("H134").Value = 4321.7144
$ws.Range("J134").Value = 3699
$ws.Range("L134").Value = 11097
$ws.Range("N134").Value = -16167
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 773330.3
$ws.Range("I136").Value = 913072.2
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 2739216.6
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -2736666.6
$ws.Range("N136").Value = -19350

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 107: H107,J107,L107,N107
$ws.Range("H107").Value = 1127.125
$ws.Range("J107").Value = 1163.6
$ws.Range("L107").Value = 3490.8
$ws.Range("N107").Value = -7330.799999999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46: H46
$ws.Range("H46").Value = 11324.5
# Row 122: H122,I122,K122,M122
$ws.Range("H122").Value = 4590.8335
$ws.Range("I122").Value = 2219.8
$ws.Range("K122").Value = 6659.400000000001
$ws.Range("M122").Value = -4209.400000000001
# Row 126: H126,I126,K126,M126
$ws.Range("H126").Value = 3826.575
$ws.Range("I126").Value = 3023.7083
$ws.Range("K126").Value = 9071.124899999999
$ws.Range("M126").Value = -6601.124899999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22: H22,I22,K22,M22
$ws.Range("H22").Value = 1194.2
$ws.Range("I22").Value = 991
$ws.Range("K22").Value = 991
$ws.Range("M22").Value = -696
# Row 27: H27,I27,K27,M27
$ws.Range("H27").Value = 1194.2
$ws.Range("I27").Value = 991
$ws.Range("K27").Value = 991
$ws.Range("M27").Value = -884
# Row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 4428.143
$ws.Range("I61").Value = 4499
$ws.Range("J61").Value = 4333.6665
$ws.Range("K61").Value = 4499
$ws.Range("L61").Value = 4333.6665
$ws.Range("M61").Value = -4297
$ws.Range("N61").Value = -4737.6665
# Row 113: H113,I113,J113,K113,L113,M113,N113
$ws.Range("H113").Value = 4428.143
$ws.Range("I113").Value = 4499
$ws.Range("J113").Value = 4333.6665
$ws.Range("K113").Value = 4499
$ws.Range("L113").Value = 4333.6665
$ws.Range("M113").Value = -2329
$ws.Range("N113").Value = -8673.666499999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100: H100,I100,K100,M100
$ws.Range("H100").Value = 910614
$ws.Range("I100").Value = 1334242.8
$ws.Range("K100").Value = 2668485.6
$ws.Range("M100").Value = -2667944.6
# Row 115: H115,J115,L115,N115
$ws.Range("H115").Value = 79996
$ws.Range("J115").Value = 79996
$ws.Range("L115").Value = 79996
$ws.Range("N115").Value = -83130
# Row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 4103.447
$ws.Range("I122").Value = 3637.641
$ws.Range("J122").Value = 6374.25
$ws.Range("K122").Value = 10912.923
$ws.Range("L122").Value = 19122.75
$ws.Range("M122").Value = -8462.923000000001
$ws.Range("N122").Value = -24022.75
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 50507104
$ws.Range("I132").Value = 9260229
$ws.Range("K132").Value = 27780687
$ws.Range("M132").Value = -27778157

Write-Output "Updated 230 cells across 8 sheets"
